$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column F (dSF) for rows 2-33, per the diff.
$newValues = @{
    2  = -3
    3  = 1
    4  = -1
    5  = -1
    6  = -2
    7  = 2
    8  = -1
    9  = 5
    10 = 0
    11 = -4
    12 = 0
    13 = 1
    14 = 3
    15 = -2
    16 = 2
    17 = 4
    18 = 3
    19 = 1
    20 = 6
    21 = 2
    22 = 0
    23 = 5
    24 = 0
    25 = 1
    26 = -2
    27 = 1
    28 = 1
    29 = 2
    30 = 3
    31 = -2
    32 = -2
    33 = -5
}

foreach ($row in $newValues.Keys) {
    $ws.Range("F$row").Value = $newValues[$row]
}
